# XLSX import : ignore unknown columns
# Rename the "Site web (URL)" header (column E, row 1) to "Site web"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E1").Value = "Site web"

# Update the active cell selection to match the author's final cursor position
$ws.Range("F15").Select() | Out-Null
